# Add two new Mac-Address rows (31 and 32) to the
# master-reg_center_user_machine worksheet, mirroring the existing rows'
# shape (regcntr_id, usr_id, machine_id, lang_code, is_active, cr_by, cr_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 31; RegCntrId = 10001; UsrId = 110030; MachineId = 10030 },
    @{ Row = 32; RegCntrId = 10001; UsrId = 110031; MachineId = 10031 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.RegCntrId
    $ws.Cells.Item($row, 2).Value = $r.UsrId
    $ws.Cells.Item($row, 3).Value = $r.MachineId
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Match the saved view state from the source workbook (scrolled down,
# single cell selected).
$ws.Range("C29").Select()
